# "Banner as 1 image." -- add a new color-swatch label pair (Body Head /
# Side Head) to the slide, and refresh the cached datetimeFigureOut
# placeholder text (slide-master + every slide layout) from 3/9/11 to
# 11/18/11.

$p = $ppt.ActivePresentation

# ---------------------------------------------------------------------
# 1) Refresh the cached "datetimeFigureOut" field text everywhere it
#    appears: the slide master and all slide layouts.
# ---------------------------------------------------------------------
function Update-DateFields($oomObj) {
    $shapes = $oomObj.Shapes
    for ($i = 1; $i -le $shapes.Count; $i++) {
        $shp = $shapes.Item($i)
        if ($shp.HasTextFrame) {
            $tr = $shp.TextFrame.TextRange
            if ($tr.Text -eq "3/9/11") {
                $tr.Text = "11/18/11"
            }
        }
    }
}

Update-DateFields $p.SlideMaster

for ($li = 1; $li -le $p.SlideMaster.CustomLayouts.Count; $li++) {
    Update-DateFields $p.SlideMaster.CustomLayouts.Item($li)
}

# ---------------------------------------------------------------------
# 2) Add the new "Body Head" / "Side Head" textboxes to slide 1 -- one
#    more color-swatch label pair, styled like the existing ones.
#    Duplicating the existing labels (instead of Shapes.AddTextbox)
#    keeps every OOXML detail (noFill, lstStyle, rtlCol, spAutoFit,
#    dirty/smtClean flags, endParaRPr) identical to authoring by hand.
# ---------------------------------------------------------------------
$s = $p.Slides.Item(1)

$bodyHeadSrc = $null
$sideHeadSrc = $null
for ($i = 1; $i -le $s.Shapes.Count; $i++) {
    $sh = $s.Shapes.Item($i)
    if ($sh.Name -eq "TextBox 4") { $bodyHeadSrc = $sh }
    if ($sh.Name -eq "TextBox 7") { $sideHeadSrc = $sh }
}

$emuPerPt = 914400 / 72
$fudge = 0.00001

# -- new "Body Head" label (id 11 / "TextBox 10") --
$bodyDup = $bodyHeadSrc.Duplicate()
$newBody = $bodyDup.Item(1)
$newBody.Name = "TextBox 10"
$newBody.Left = 5502714 / $emuPerPt + $fudge
$newBody.Top = 1429869 / $emuPerPt + $fudge
$newBody.Width = 1539003 / $emuPerPt + $fudge
$newBody.Height = 461665 / $emuPerPt + $fudge
$newBody.TextFrame.TextRange.Font.Color.RGB = 34504   # C88600

# -- new "Side Head" label (id 12 / "TextBox 11") --
$sideDup = $sideHeadSrc.Duplicate()
$newSide = $sideDup.Item(1)
$newSide.Name = "TextBox 11"
$newSide.Left = 550227 / $emuPerPt + $fudge
$newSide.Top = 1846165 / $emuPerPt + $fudge
$newSide.Width = 1226743 / $emuPerPt + $fudge
$newSide.Height = 400110 / $emuPerPt + $fudge
$newSide.TextFrame.TextRange.Font.Color.RGB = 5816822  # F6C158

Write-Host "Added shapes: $($newBody.Name) (id=$($newBody.Id)), $($newSide.Name) (id=$($newSide.Id))"
